# Apply updated crafting-profit figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 54995
$ws.Range("J57").Value = 54995
$ws.Range("L57").Value = 164985
$ws.Range("N57").Value = -165983

$ws.Range("H112").Value = 1010.88
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 1016.1739
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 3048.5217
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -5264.5217

$ws.Range("H132").Value = 5213039
$ws.Range("I132").Value = 6255489
$ws.Range("J132").Value = 788.25
$ws.Range("K132").Value = 18766467
$ws.Range("L132").Value = 2364.75
$ws.Range("M132").Value = -18763937
$ws.Range("N132").Value = -7424.75

$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

$ws.Range("H135").Value = 505.13333
$ws.Range("I135").Value = 505.13333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4546.19997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2011.19997
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 1301.0944
$ws.Range("I137").Value = 1012.5
$ws.Range("J137").Value = 1476
$ws.Range("K137").Value = 3037.5
$ws.Range("L137").Value = 4428
$ws.Range("M137").Value = -487.5
$ws.Range("N137").Value = -9528

$ws.Range("H138").Value = 3783.1143
$ws.Range("I138").Value = 1757.1724
$ws.Range("J138").Value = 5216.0977
$ws.Range("K138").Value = 5271.5172
$ws.Range("L138").Value = 15648.2931
$ws.Range("M138").Value = -131.5172000000002
$ws.Range("N138").Value = -25928.2931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28343.676
$ws.Range("I2").Value = 1034.75
$ws.Range("K2").Value = 1034.75
$ws.Range("M2").Value = -921.75

$ws.Range("H45").Value = 1518.2222
$ws.Range("J45").Value = 2184.625
$ws.Range("L45").Value = 2184.625
$ws.Range("N45").Value = -2938.625

$ws.Range("H61").Value = 1905.1111
$ws.Range("I61").Value = 1004.5
$ws.Range("J61").Value = 2625.6
$ws.Range("K61").Value = 1004.5
$ws.Range("L61").Value = 2625.6
$ws.Range("M61").Value = -792.5
$ws.Range("N61").Value = -3049.6

$ws.Range("H74").Value = 2610.6365
$ws.Range("I74").Value = 1369.1428
$ws.Range("J74").Value = 3190
$ws.Range("K74").Value = 1369.1428
$ws.Range("L74").Value = 3190
$ws.Range("M74").Value = -495.1428000000001
$ws.Range("N74").Value = -4938

$ws.Range("H77").Value = 2610.6365
$ws.Range("I77").Value = 1369.1428
$ws.Range("J77").Value = 3190
$ws.Range("K77").Value = 6845.714
$ws.Range("L77").Value = 15950
$ws.Range("M77").Value = -2477.714
$ws.Range("N77").Value = -24686

$ws.Range("H116").Value = 28343.676
$ws.Range("I116").Value = 1034.75
$ws.Range("K116").Value = 1034.75
$ws.Range("M116").Value = 1259.25

$ws.Range("H136").Value = 1905.1111
$ws.Range("I136").Value = 1004.5
$ws.Range("J136").Value = 2625.6
$ws.Range("K136").Value = 3013.5
$ws.Range("L136").Value = 7876.799999999999
$ws.Range("M136").Value = -463.5
$ws.Range("N136").Value = -12976.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28343.676
$ws.Range("I3").Value = 1034.75
$ws.Range("K3").Value = 1034.75
$ws.Range("M3").Value = -920.75

$ws.Range("H105").Value = 75975.19
$ws.Range("I105").Value = 51803.45
$ws.Range("J105").Value = 145037.28
$ws.Range("K105").Value = 51803.45
$ws.Range("L105").Value = 145037.28
$ws.Range("M105").Value = -50056.45
$ws.Range("N105").Value = -148531.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19249.31
$ws.Range("J31").Value = 2411.0527
$ws.Range("L31").Value = 2411.0527
$ws.Range("N31").Value = -3001.0527

$ws.Range("H34").Value = 19249.31
$ws.Range("J34").Value = 2411.0527
$ws.Range("L34").Value = 2411.0527
$ws.Range("N34").Value = -2815.0527

$ws.Range("H52").Value = 45513.332
$ws.Range("J52").Value = 45513.332
$ws.Range("L52").Value = 45513.332
$ws.Range("N52").Value = -46101.332

$ws.Range("H58").Value = 10396.448
$ws.Range("I58").Value = 1729.2632
$ws.Range("J58").Value = 26864.1
$ws.Range("K58").Value = 1729.2632
$ws.Range("L58").Value = 26864.1
$ws.Range("M58").Value = -1526.2632
$ws.Range("N58").Value = -27270.1

$ws.Range("H132").Value = 2884.4
$ws.Range("I132").Value = 2911.6667
$ws.Range("K132").Value = 8735.000100000001
$ws.Range("M132").Value = -6205.000100000001

$ws.Range("H136").Value = 10396.448
$ws.Range("I136").Value = 1729.2632
$ws.Range("J136").Value = 26864.1
$ws.Range("K136").Value = 5187.7896
$ws.Range("L136").Value = 80592.29999999999
$ws.Range("M136").Value = -2637.7896
$ws.Range("N136").Value = -85692.29999999999

$ws.Range("H139").Value = 33263
$ws.Range("J139").Value = 34540
$ws.Range("L139").Value = 34540
$ws.Range("N139").Value = -44820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 367.25
$ws.Range("I92").Value = 100
$ws.Range("J92").Value = 456.33334
$ws.Range("K92").Value = 300
$ws.Range("L92").Value = 1369.00002
$ws.Range("M92").Value = 948
$ws.Range("N92").Value = -3865.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 431983.56
$ws.Range("I102").Value = 5829.5
$ws.Range("K102").Value = 5829.5
$ws.Range("M102").Value = -4207.5

$ws.Range("H137").Value = 68000
$ws.Range("J137").Value = 68000
$ws.Range("L137").Value = 68000
$ws.Range("N137").Value = -78200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18150.666
$ws.Range("I132").Value = 18150.666
$ws.Range("K132").Value = 54451.99800000001
$ws.Range("M132").Value = -51921.99800000001

$ws.Range("H136").Value = 3436.7273
$ws.Range("I136").Value = 2875.5
$ws.Range("K136").Value = 8626.5
$ws.Range("M136").Value = -6076.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2417864.5
$ws.Range("J62").Value = 2600
$ws.Range("L62").Value = 2600
$ws.Range("N62").Value = -3848

$ws.Range("H65").Value = 2417864.5
$ws.Range("J65").Value = 2600
$ws.Range("L65").Value = 13000
$ws.Range("N65").Value = -19240

$ws.Range("H132").Value = 11621.667
$ws.Range("I132").Value = 12946.2
$ws.Range("K132").Value = 38838.60000000001
$ws.Range("M132").Value = -36308.60000000001

$ws.Range("H136").Value = 1526.8182
$ws.Range("I136").Value = 999.2857
$ws.Range("K136").Value = 2997.8571
$ws.Range("M136").Value = -447.8571000000002

$ws.Range("H139").Value = 65635
$ws.Range("J139").Value = 65635
$ws.Range("L139").Value = 65635
$ws.Range("N139").Value = -75915
